# Weekly update: insert a new price record for "Cebollín" (Femacal de La
# Calera) ahead of the existing rows, shifting the rest of the table down
# by one row (dimension grows from A1:R492 to A1:R493).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row above the current row 384; Excel shifts rows
# 384:492 down to 385:493 and copies row 384's formatting (incl. the
# date number-format in column D) into the new blank row.
$ws.Rows("384:384").Insert()

$ws.Cells.Item(384, 1).Value2 = 3
$ws.Cells.Item(384, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(384, 3).Value2 = "Coquimbo"
$ws.Cells.Item(384, 4).Value2 = 44736
$ws.Cells.Item(384, 5).Value2 = 5
$ws.Cells.Item(384, 6).Value2 = 100112037
$ws.Cells.Item(384, 7).Value2 = "Cebollín"
$ws.Cells.Item(384, 8).Value2 = "Sin especificar"
$ws.Cells.Item(384, 9).Value2 = "Primera"
$ws.Cells.Item(384, 10).Value2 = 250
$ws.Cells.Item(384, 11).Value2 = 6500
$ws.Cells.Item(384, 12).Value2 = 7000
$ws.Cells.Item(384, 13).Value2 = 6760
$ws.Cells.Item(384, 14).Value2 = "$/paquete 36 unidades"
$ws.Cells.Item(384, 15).Value2 = "Provincia de Quillota"
$ws.Cells.Item(384, 16).Value2 = 188
$ws.Cells.Item(384, 17).Value2 = 36
$ws.Cells.Item(384, 18).Value2 = "Hortaliza"
